$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9112048919340322
$ws.Range("C2").Value = 1.133786848072562
$ws.Range("D2").Value = 0.1126721327937182
$ws.Range("E2").Value = 0.9070294784580499
$ws.Range("F2").Value = 1.139435185830596
$ws.Range("G2").Value = 0.1133786848072562

$ws.Range("B3").Value = -0.4402577480568004
$ws.Range("C3").Value = -0.36281179138322
$ws.Range("D3").Value = 0.1028323251606915
$ws.Range("E3").Value = -0.4308390022675737
$ws.Range("F3").Value = -0.3662529298129745
$ws.Range("G3").Value = 0.1133786848072562

$ws.Range("B4").Value = -2.031350183375844
$ws.Range("C4").Value = -2.448979591836735
$ws.Range("D4").Value = -1.253243669098878
$ws.Range("E4").Value = -2.040816326530612
$ws.Range("F4").Value = -2.444198566859693
$ws.Range("G4").Value = -1.26984126984127

$ws.Range("B5").Value = -0.002692718193406324
$ws.Range("C5").Value = -0.2040816326530612
$ws.Range("D5").Value = -0.1855199464135275
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = -0.1878768028333454
$ws.Range("G5").Value = -0.18140589569161

$ws.Range("B6").Value = 1.541215953830037
$ws.Range("C6").Value = 1.927437641723356
$ws.Range("D6").Value = 0.1231218379758774
$ws.Range("E6").Value = 1.541950113378685
$ws.Range("F6").Value = 1.945556573671575
$ws.Range("G6").Value = 0.1133786848072562

$ws.Range("B7").Value = -0.3802581138671877
$ws.Range("C7").Value = 0.18140589569161
$ws.Range("D7").Value = 0.7815463054709451
$ws.Range("E7").Value = -0.3854875283446712
$ws.Range("F7").Value = 0.1838007771989899
$ws.Range("G7").Value = 0.7709750566893424

$ws.Range("B8").Value = 0.5145133491967829
$ws.Range("C8").Value = 0.453514739229025
$ws.Range("D8").Value = -0.04511345882168556
$ws.Range("E8").Value = 0.5215419501133787
$ws.Range("F8").Value = 0.449486283507921
$ws.Range("G8").Value = -0.04535147392290249

$ws.Range("B9").Value = 0.4907005132262501
$ws.Range("C9").Value = 1.224489795918367
$ws.Range("D9").Value = 0.4794669711284663
$ws.Range("E9").Value = 0.4988662131519275
$ws.Range("F9").Value = 1.220539079278846
$ws.Range("G9").Value = 0.4761904761904762

$ws.Range("B10").Value = 1.165374245474661
$ws.Range("C10").Value = 2.244897959183673
$ws.Range("D10").Value = 0.3733966858293393
$ws.Range("E10").Value = 1.179138321995465
$ws.Range("F10").Value = 2.245223706833127
$ws.Range("G10").Value = 0.3854875283446712

$ws.Range("B11").Value = 0.8777519492556363
$ws.Range("C11").Value = 0.4308390022675737
$ws.Range("D11").Value = -0.2530379391380345
$ws.Range("E11").Value = 0.8843537414965986
$ws.Range("F11").Value = 0.4345807984913791
$ws.Range("G11").Value = -0.2494331065759637
